$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.523.35'
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('D3').Value = '1.820.48'
$ws.Range('E3').Value = '  -0.41%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.85'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.24%  '
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5150'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.62%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3883'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.65%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08452'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +8.66%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.82'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.109'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.04%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.439'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.64%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '21.00'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.59%  '
$ws.Range('E14').Value = '  +0.12%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.515'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.97%  '
$ws.Range('D16').Value = '1.824.44'
$ws.Range('E16').Value = '  -0.21%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001139'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.34%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '92.82'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06689'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.75'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.44%  '
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.090'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').Value = '28.555.92'
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.41'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.74%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.274'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '21.12'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.29%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '158.63'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.88%  '
$ws.Range('D28').Value = '2.029.74'
$ws.Range('E28').Value = '  -0.98%  '
$ws.Range('E29').Value = '  +0.20%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '126.01'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.28%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1085'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.34%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.096'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.88%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.740'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.11%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.07546'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.78%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.679'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.75%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.2231'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.87%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02365'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.34%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.192'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.51%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.735'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6326'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.38%  '
$ws.Range('E41').Value = '  -1.49%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.195'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.28%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.401'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.48%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.61'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.777'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.80%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5943'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '125.94'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.03%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.989'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.57%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.200'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.34%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06979'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.33%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '74.43'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.26%  '
